$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")
$ws.Activate()

# Replace the "Good Morning" greeting text with "GIT UPDATE"
$ws.Range("E8").Value = "GIT UPDATE"

# Record the new active cell selection on the sheet (matches the authored commit)
$ws.Range("E8").Select()
